$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the D/J/K/L/M/P values among rows 2, 3 and 5:
#   new row2 = old row5
#   new row3 = old row2
#   new row5 = old row3
$cols = @("D", "J", "K", "L", "M", "P")

$old2 = @{}
$old3 = @{}
$old5 = @{}
foreach ($col in $cols) {
    $old2[$col] = $ws.Range($col + "2").Value2
    $old3[$col] = $ws.Range($col + "3").Value2
    $old5[$col] = $ws.Range($col + "5").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $old5[$col]
    $ws.Range($col + "3").Value2 = $old2[$col]
    $ws.Range($col + "5").Value2 = $old3[$col]
}
